$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Copy cell formatting (styles) to the cells that need a new
#    style BEFORE their current source cells get overwritten.
# -----------------------------------------------------------------
# B4 and (temporarily) D4 need the style currently on C4.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# E4 needs the style currently on D4.
$ws.Range("D4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

# D4 (new) needs the style C4 used to have.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null

# C4 (new) needs the header-row style (same as A3/B3/.../H3).
$ws.Range("A3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# -----------------------------------------------------------------
# 2) Header row (row 3). Order matters: new shared strings get
#    appended to the table in the order they are first written, so
#    write them in the same order the target workbook expects
#    (#, type, publication title).
# -----------------------------------------------------------------
$ws.Range("A3").Value = "№"
$ws.Range("C3").Value = "тип"
$ws.Range("B3").Value = "Назва публікації"
$ws.Range("D3").Value = "Вихідні дані"
$ws.Range("E3").Value = "дата видання"

# -----------------------------------------------------------------
# 3) Data row (row 4): shift "Вихідні дані"/"дата видання" values
#    one column to the right, clear the old ПІБ / "Економіка..." cells,
#    and move "Модернизация..." into column B.
# -----------------------------------------------------------------
$ws.Range("A4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("B4").Value = "Модернизация: путь в будущее"
$ws.Range("D4").Value = "Економіка розвитку, 2015.  -№3. - С. 21-36. "
$ws.Range("E4").Value = 42112

# -----------------------------------------------------------------
# 4) Column widths.
# -----------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 5.666666666666667
$ws.Columns("B").ColumnWidth = 20.166666666666668
$ws.Columns("C").ColumnWidth = 17.5
$ws.Columns("D").ColumnWidth = 20.333333333333336
$ws.Columns("E").ColumnWidth = 14.666666666666666
$ws.Columns("H").ColumnWidth = 19.666666666666668

# -----------------------------------------------------------------
# 5) Selection.
# -----------------------------------------------------------------
$ws.Range("C5").Select() | Out-Null
